$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169576644897461
$ws.Range("B1").Value = 2.439228534698486
$ws.Range("D1").Value = 2.366158962249756
$ws.Range("E1").Value = 1.233809113502502
